$wb = $excel.ActiveWorkbook

# Updates apply identically to the "展览" and "全部类型" sheets, which mirror
# the same event data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F5").Value = 2782
    $ws.Range("F9").Value = 1486
    $ws.Range("F11").Value = 69
    $ws.Range("F13").Value = 1240
    $ws.Range("F14").Value = 9
    $ws.Range("F15").Value = 383
    $ws.Range("F16").Value = 332
    $ws.Range("F17").Value = 51
    $ws.Range("F22").Value = 2733
    $ws.Range("F23").Value = 328
    $ws.Range("G23").Value = 55
    $ws.Range("F24").Value = 5
}
